$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (columns B..Q) - identical for every data row (2..26)
$values = @(
    0.6408044419705359,
    -265.1326953808737,
    0.4275118738973983,
    0.7064707130608152,
    0.7236178716181977,
    0.2132339996864685,
    157.9878643119225,
    0.1594276330412308,
    0.09697370464141555,
    0.1282006688413232,
    0.2405142646481177,
    0.4617726710043249,
    0.2163006006629874,
    0.4814312896101858,
    29.09073025240775,
    44.93611597569436
)

for ($row = 2; $row -le 26; $row++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # Column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
